# Update header labels so that Power BI can auto-detect the first row as a
# table header. Prefix the year labels with "Ano " (or "Intervalo " for the
# sheet that uses period ranges) on each sheet's first row.

$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3, 5: header row B1:E1 = 2015 / 2030 / 2040 / 2050 -> "Ano <year>"
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet 4: header row B1:E1 = 2015 / 2015-2030 / 2031-2040 / 2041-2050 -> "Intervalo <period>"
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("B1").Value = "Intervalo 2015"
$ws4.Range("C1").Value = "Intervalo 2015-2030"
$ws4.Range("D1").Value = "Intervalo 2031-2040"
$ws4.Range("E1").Value = "Intervalo 2041-2050"

# Sheet 6: only B1 = 2015 -> "Ano 2015"
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Range("B1").Value = "Ano 2015"
